{"js": "// Logbook update: append two new dated entries at the end of the document,\n// right before the final trailing blank paragraph, mirroring the existing\n// \"date paragraph\" / \"content paragraph\" / \"blank separator\" pattern used\n// throughout the rest of the log.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The document always ends with a blank paragraph (the separator after the\n// most recent entry). Insert the new entries immediately before it so the\n// trailing blank paragraph stays last, exactly like every earlier entry.\nconst trailingBlankParagraph = paragraphs.items[paragraphs.items.length - 1];\n\nconst newContent = [\n  \"\",\n  \"07/05/23\",\n  \"Finished GUI elements of Tkinter implementation.\",\n  \"\",\n  \"12/05/23\",\n  \"Finished decision logic of Tkinter implementation, as well as flask implementation. Submitted.\"\n];\n\nfor (const text of newContent) {\n  trailingBlankParagraph.insertParagraph(text, Word.InsertLocation.before);\n}\n\nawait context.sync();\n", "ps1": "# Logbook update: append two new dated entries at the end of the document,\n# right before the final trailing blank paragraph, mirroring the existing\n# \"date paragraph\" / \"content paragraph\" / \"blank separator\" pattern used\n# throughout the rest of the log.\n\n$d = $word.ActiveDocument\n\n# The document always ends with a blank paragraph (the separator after the\n# most recent entry). Position the insertion point at the very start of\n# that paragraph so the new entries land before it and it stays last,\n# exactly like every earlier entry's trailing blank line.\n$lastParagraph = $d.Paragraphs.Last\n$insertionPoint = $lastParagraph.Range\n$insertionPoint.Collapse(1)  # wdCollapseStart\n\n$newLines = @(\n  \"\",\n  \"07/05/23\",\n  \"Finished GUI elements of Tkinter implementation.\",\n  \"\",\n  \"12/05/23\",\n  \"Finished decision logic of Tkinter implementation, as well as flask implementation. Submitted.\"\n)\n\nforeach ($line in $newLines) {\n  $insertionPoint.InsertBefore($line + \"`r\")\n  $insertionPoint.Collapse(0)  # wdCollapseEnd -- advance past the line just inserted\n}\n"}
